$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the region code from DE to AT
$ws.Range("B2").Value = "AT"

# Update the outside temperature value from -10 to -12
$ws.Range("D2").Value = -12

# Move the active selection to B3, matching the post-edit selection state
$ws.Range("B3").Select()
